$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set D column (Price) cells - force text format to avoid numeric auto-conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.903.51'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.766.97'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '328.52'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4529'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3523'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '41.95'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07377'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.092'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.70'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.004'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.179'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.770.28'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.46'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001060'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06443'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.93'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.754'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.933.69'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.19'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.49'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.13'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.972.82'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.93'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.076'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09179'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.607'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.666'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.82'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02280'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06110'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.942'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6246'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.178'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.781'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.27'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.736'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5848'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.52'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.928'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.129'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06828'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.143'

# Set E column (Volume 1h) cells
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("E3").Value = '  +0.66%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("E6").Value = '  -0.05%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  -0.90%  '
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("E11").Value = '  +0.68%  '
$ws.Range("E12").Value = '  -0.05%  '
$ws.Range("E13").Value = '  -0.19%  '
$ws.Range("E14").Value = '  +0.42%  '
$ws.Range("E15").Value = '  +0.41%  '
$ws.Range("E16").Value = '  +1.19%  '
$ws.Range("E17").Value = '  -1.03%  '
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("E21").Value = '  -0.67%  '
$ws.Range("E23").Value = '  +0.95%  '
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("E26").Value = '  -3.91%  '
$ws.Range("E27").Value = '  -0.23%  '
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("E29").Value = '  +2.81%  '
$ws.Range("E30").Value = '  -1.03%  '
$ws.Range("E31").Value = '  -0.74%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("E34").Value = '  +0.32%  '
$ws.Range("E35").Value = '  +1.18%  '
$ws.Range("E36").Value = '  -0.41%  '
$ws.Range("E37").Value = '  +1.33%  '
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("E40").Value = '  -0.75%  '
$ws.Range("E41").Value = '  -0.33%  '
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("E44").Value = '  +0.76%  '
$ws.Range("E45").Value = '  +0.64%  '
$ws.Range("E47").Value = '  -0.20%  '
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("E49").Value = '  +0.15%  '
$ws.Range("E50").Value = '  -1.02%  '
$ws.Range("E51").Value = '  +1.66%  '
